$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top, pushing all data down by one row.
$ws.Rows.Item(1).Insert()

# The old stray cell (previously D7, now shifted to D8) is no longer needed;
# remove that now-empty row so the used range ends at row 7 again.
$ws.Rows.Item(8).Delete()

# Update the active selection to match the new state.
$ws.Range("A9").Select() | Out-Null
